$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.564.47'
$ws.Range('E2').Value = '  +1.47%  '

# Row 3
$ws.Range('D3').Value = '1.884.44'
$ws.Range('E3').Value = '  +1.50%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.26%  '

# Row 6
$ws.Range('E6').Value = '  -0.04%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4764'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.62%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2918'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.13%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06536'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.33%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.08'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.67%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '98.04'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.10%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07723'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.36%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7422'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.47%  '

# Row 14
$ws.Range('D14').Value = '1.881.91'
$ws.Range('E14').Value = '  +1.21%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.156'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.99%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '274.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.93%  '

# Row 17
$ws.Range('D17').Value = '30.553.07'
$ws.Range('E17').Value = '  +1.51%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.09%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007578'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.50%  '

# Row 21
$ws.Range('D21').Value = '2.128.53'
$ws.Range('E21').Value = '  +1.11%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.06%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.261'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.39%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.203'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.81%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.340'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.92%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.99%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.82%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.952'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.99%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1006'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.49%  '

# Row 30
$ws.Range('E30').Value = '  +0.14%  '

# Row 31
$ws.Range('E31').Value = '  +4.61%  '

# Row 32
$ws.Range('E32').Value = '  +3.15%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.131'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.04%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04816'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.68%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.134'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.02%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7027'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.85%  '

# Row 37
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.03%  '

# Row 38
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.715'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.07%  '

# Row 39
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01870'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.64%  '

# Row 40
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.749'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.17%  '

# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.332'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.01%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.72%  '

# Row 43
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.77'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.06%  '

# Row 44
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4231'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.92%  '

# Row 45
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8417'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.45%  '

# Row 46
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.02%  '

# Row 47
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.91'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.58%  '

# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.295'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.66%  '

# Row 49
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.119'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.26%  '

# Row 50
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.69'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.94%  '

# Row 51
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '917.84'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.52%  '
